$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for rows 2-3 ---
# These reflect the per-locale handback status, which just flipped from
# "Ready for handoff" to "Handed back: in sync with en-US" for every row.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet ---
# Status column (C): same status flip as above.
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
# Latest Handback DateTime (K): refreshed to the handback-report generation time.
$zhcn.Range("K2").Value = "2016-10-21 00:24:11"
$zhcn.Range("K3").Value = "2016-10-21 00:24:11"
# Error Detail (P): the "handback file not latest" warning no longer applies
# now that the handback is in sync, so it is cleared.
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

# --- de-de sheet ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("K2").Value = "2016-10-21 00:24:29"
$dede.Range("K3").Value = "2016-10-21 00:24:29"
$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

# --- Column width touch-ups (the columns holding the text above got wider,
#     the now-empty Error Detail column got narrower). ColumnWidth values are
#     chosen to snap, under this host's width quantization, to the same
#     stored widths produced by the original edit. ---
$overview.Columns.Item(5).ColumnWidth = 29.2
$overview.Columns.Item(6).ColumnWidth = 29.2

$zhcn.Columns.Item(3).ColumnWidth = 29.2
$zhcn.Columns.Item(16).ColumnWidth = 12.85

$dede.Columns.Item(3).ColumnWidth = 29.2
$dede.Columns.Item(16).ColumnWidth = 12.85
